# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.993.58'
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('D3').Value = '3.522.46'
$ws.Range('E3').Value = '  -2.32%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '602.76'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '143.62'
$ws.Range('E6').Value = '  -3.48%  '
$ws.Range('D7').Value = '3.523.74'
$ws.Range('E7').Value = '  -2.32%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.512'
$ws.Range('E9').Value = '  +4.66%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '7.81'
$ws.Range('E10').Value = '  -3.04%  '
$ws.Range('E11').Value = '  -4.66%  '
$ws.Range('E12').Value = '  -2.10%  '
$ws.Range('D13').Value = '4.130.70'
$ws.Range('E13').Value = '  -1.98%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.0000193'
$ws.Range('E14').Value = '  -8.07%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '28.31'
$ws.Range('E15').Value = '  -5.95%  '
$ws.Range('D16').Value = '3.530.28'
$ws.Range('E16').Value = '  -2.01%  '
$ws.Range('E17').Value = '  +0.94%  '
$ws.Range('D18').Value = '65.912.83'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.98'
$ws.Range('E19').Value = '  -4.65%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.17'
$ws.Range('E20').Value = '  -2.91%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.55'
$ws.Range('E21').Value = '  -3.52%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '421.06'
$ws.Range('E22').Value = '  -1.91%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.593'
$ws.Range('E23').Value = '  -4.70%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '76.88'
$ws.Range('E24').Value = '  -2.92%  '
$ws.Range('D25').Value = '3.672.48'
$ws.Range('E25').Value = '  -1.98%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -7.16%  '
$ws.Range('E28').Value = '  -2.58%  '
$ws.Range('E29').Value = '  -6.12%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.84'
$ws.Range('E30').Value = '  -5.03%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').Value = '3.537.73'
$ws.Range('E32').Value = '  -1.84%  '
$ws.Range('E33').Value = '  -0.65%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '24.13'
$ws.Range('E34').Value = '  -5.61%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.34'
$ws.Range('E36').Value = '  -8.32%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '7.55'
$ws.Range('E37').Value = '  -4.14%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '178.12'
$ws.Range('E38').Value = '  +0.73%  '
$ws.Range('E39').Value = '  -5.25%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.20'
$ws.Range('E40').Value = '  -7.77%  '
$ws.Range('E41').Value = '  -5.16%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '4.96'
$ws.Range('E42').Value = '  -5.55%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.857'
$ws.Range('E43').Value = '  -4.77%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '45.33'
$ws.Range('E44').Value = '  -1.95%  '
$ws.Range('E45').Value = '  -8.59%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.39'
$ws.Range('E47').Value = '  -7.88%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '23.72'
$ws.Range('E48').Value = '  -2.81%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '7.03'
$ws.Range('E49').Value = '  -2.40%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.11'
$ws.Range('E50').Value = '  -7.07%  '
$ws.Range('E51').Value = '  -5.29%  '
